$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") from 45188 to 45189 for all data rows (2..527)
$lastDataRow = 527
for ($r = 2; $r -le $lastDataRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45188) {
        $cCell.Value2 = 45189
    }
}

# 2) Ensure row 527 gets an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(527).RowHeight = 15

# 3) Append two new rows (528 and 529) with the new entries
$newRows = @(
    @{ Row = 528; A = "A 43765-2023"; B = 45187; C = 45189; D = "NORRBOTTENS LÄN"; E = "LULEÅ"; G = 13.7 },
    @{ Row = 529; A = "A 43861-2023"; B = 45187; C = 45189; D = "NORRBOTTENS LÄN"; E = "LULEÅ"; G = 5.1 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    $ws.Cells.Item($r, 1).Value2 = $nr.A

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value2 = $nr.B
    $bCell.NumberFormat = "YYYY-MM-DD"

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value2 = $nr.C
    $cCell.NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value2 = $nr.D
    $ws.Cells.Item($r, 5).Value2 = $nr.E

    $ws.Cells.Item($r, 7).Value2 = $nr.G

    # Columns H..Q are all 0
    for ($col = 8; $col -le 17; $col++) {
        $ws.Cells.Item($r, $col).Value2 = 0
    }

    # Column R keeps the wrap-text style, no content
    $rCell = $ws.Cells.Item($r, 18)
    $rCell.WrapText = $true
}

# Row 528 gets an explicit row height (ht="15" customHeight="1"), row 529 keeps the default
$ws.Rows.Item(528).RowHeight = 15
